# The sheet holds one job posting per row (A=index, B=_id, C=fulltext).
# Row 21 ("639134d6...") was dropped from the source DB. The updated export
# re-wrote only the B (_id) and C (fulltext) columns by shifting every
# subsequent record up by one position, leaving column A's running index
# untouched, and then trimmed the now-duplicate trailing row (37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the _id/fulltext pairs of rows 22:37 up into rows 21:36.
$src = $ws.Range("B22:C37")
$dst = $ws.Range("B21:C36")
$src.Copy($dst)

# Remove the now-redundant last row so the sheet shrinks from 37 to 36 rows.
$ws.Rows(37).Delete()
